$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2144846796657382
$ws.Range("C2").Value = 0.5264623955431755
$ws.Range("J2").Value = 0.01114206128133705
$ws.Range("P2").Value = 0.1448467966573816
$ws.Range("S2").Value = 0.1030640668523677
$ws.Range("B3").Value = 0.01036269430051814
$ws.Range("C3").Value = 0.0155440414507772
$ws.Range("J3").Value = 0.03626943005181347
$ws.Range("P3").Value = 0.7461139896373057
$ws.Range("S3").Value = 0.1917098445595855
$ws.Range("J4").Value = 0.07462686567164178
$ws.Range("O4").Value = 0.01492537313432836
$ws.Range("P4").Value = 0.6417910447761194
$ws.Range("S4").Value = 0.2686567164179104
$ws.Range("B6").Value = 0.1041666666666667
$ws.Range("D6").Value = 0.008333333333333333
$ws.Range("F6").Value = 0.075
$ws.Range("J6").Value = 0.2166666666666667
$ws.Range("O6").Value = 0.025
$ws.Range("Q6").Value = 0.1625
$ws.Range("R6").Value = 0.07083333333333333
$ws.Range("S6").Value = 0.3375
$ws.Range("B7").Value = 0.1692307692307692
$ws.Range("D7").Value = 0.03076923076923077
$ws.Range("E7").Value = 0.003846153846153846
$ws.Range("F7").Value = 0.05384615384615385
$ws.Range("J7").Value = 0.1
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.03846153846153846
$ws.Range("S7").Value = 0.4346153846153846
$ws.Range("B8").Value = 0.1106290672451193
$ws.Range("D8").Value = 0.02386117136659436
$ws.Range("E8").Value = 0.002169197396963124
$ws.Range("F8").Value = 0.05206073752711497
$ws.Range("J8").Value = 0.09761388286334056
$ws.Range("O8").Value = 0.03253796095444685
$ws.Range("Q8").Value = 0.1735357917570499
$ws.Range("R8").Value = 0.08459869848156182
$ws.Range("S8").Value = 0.4229934924078091
$ws.Range("B9").Value = 0.09130434782608696
$ws.Range("D9").Value = 0.03043478260869565
$ws.Range("F9").Value = 0.0782608695652174
$ws.Range("J9").Value = 0.1217391304347826
$ws.Range("O9").Value = 0.01739130434782609
$ws.Range("Q9").Value = 0.2217391304347826
$ws.Range("R9").Value = 0.06956521739130435
$ws.Range("S9").Value = 0.3695652173913043
$ws.Range("B10").Value = 0.1144624903325599
$ws.Range("D10").Value = 0.02938901778808971
$ws.Range("E10").Value = 0.0007733952049497294
$ws.Range("F10").Value = 0.07269914926527456
$ws.Range("J10").Value = 0.1160092807424594
$ws.Range("O10").Value = 0.02010827532869296
$ws.Range("Q10").Value = 0.234338747099768
$ws.Range("R10").Value = 0.05723124516627997
$ws.Range("S10").Value = 0.3549883990719258
$ws.Range("F11").Value = 0.002487562189054726
$ws.Range("G11").Value = 0.1467661691542289
$ws.Range("J11").Value = 0.09701492537313433
$ws.Range("K11").Value = 0.208955223880597
$ws.Range("L11").Value = 0.5099502487562189
$ws.Range("S11").Value = 0.03482587064676617
$ws.Range("G12").Value = 0.7633928571428571
$ws.Range("J12").Value = 0.1830357142857143
$ws.Range("K12").Value = 0.01785714285714286
$ws.Range("L12").Value = 0.008928571428571428
$ws.Range("S12").Value = 0.02678571428571428
$ws.Range("G13").Value = 0.6896551724137931
$ws.Range("J13").Value = 0.2068965517241379
$ws.Range("S13").Value = 0.103448275862069
$ws.Range("F15").Value = 0.02040816326530612
$ws.Range("H15").Value = 0.1306122448979592
$ws.Range("I15").Value = 0.08979591836734693
$ws.Range("J15").Value = 0.3102040816326531
$ws.Range("K15").Value = 0.05714285714285714
$ws.Range("M15").Value = 0.02857142857142857
$ws.Range("O15").Value = 0.07346938775510205
$ws.Range("S15").Value = 0.2897959183673469
$ws.Range("F16").Value = 0.01724137931034483
$ws.Range("H16").Value = 0.1637931034482759
$ws.Range("I16").Value = 0.07327586206896551
$ws.Range("J16").Value = 0.4008620689655172
$ws.Range("K16").Value = 0.1336206896551724
$ws.Range("M16").Value = 0.02586206896551724
$ws.Range("O16").Value = 0.06465517241379311
$ws.Range("S16").Value = 0.1206896551724138
$ws.Range("F17").Value = 0.01747572815533981
$ws.Range("H17").Value = 0.1514563106796117
$ws.Range("I17").Value = 0.116504854368932
$ws.Range("J17").Value = 0.3689320388349515
$ws.Range("K17").Value = 0.1262135922330097
$ws.Range("M17").Value = 0.01553398058252427
$ws.Range("O17").Value = 0.07378640776699029
$ws.Range("S17").Value = 0.1300970873786408
$ws.Range("F18").Value = 0.0125
$ws.Range("H18").Value = 0.20625
$ws.Range("I18").Value = 0.075
$ws.Range("J18").Value = 0.43125
$ws.Range("K18").Value = 0.09375
$ws.Range("M18").Value = 0.01875
$ws.Range("O18").Value = 0.06875000000000001
$ws.Range("S18").Value = 0.09375
$ws.Range("F19").Value = 0.02
$ws.Range("H19").Value = 0.1928571428571429
$ws.Range("I19").Value = 0.08214285714285714
$ws.Range("J19").Value = 0.3442857142857143
$ws.Range("K19").Value = 0.1278571428571429
$ws.Range("M19").Value = 0.02357142857142857
$ws.Range("O19").Value = 0.05642857142857143
$ws.Range("S19").Value = 0.1528571428571429
